# Update the South Africa Semiconductor capital structure database
# with refreshed figures (rows 2 and 3 share identical values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G"  = -0.01158682634730539
    "H"  = -0.01158682634730539
    "I"  = -0.04817250612141574
    "J"  = -0.04817250612141574
    "K"  = -4.24
    "L"  = -0.1269461077844312
    "U"  = 0.405
    "V"  = 0.03461538461538462
    "W"  = -1.675889328063241
    "X"  = 0.1691901365342467
    "Y"  = -1.845079464597488
    "Z"  = -11.92349764934419
    "AA" = 0.5743847635017191
    "AB" = 0.1188775123220451
    "AC" = 0.455507251179674
    "AD" = 7.44
    "AE" = 0.08980852227642722
    "AF" = 7.529808522276427
    "AG" = 7.124808522276427
    "AH" = 0.3915696047391036
    "AI" = 0.6441344619056345
    "AJ" = 0.3784797340087285
    "AK" = 0.6313628191574469
    "AL" = 0.667
    "AM" = 0.665
    "AN" = -9.600000000000001
    "AO" = -2.443778110944527
    "AP" = -9.193301319066359
    "AQ" = -2.451127819548872
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# O and R keep magnitude 0 but flip the signed-zero representation
# between the two rows.
$ws.Range("O2").Value = -0
$ws.Range("R2").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0

$wb.Save()
